$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer Quote")

# Fixed surplus number: the currency surcharge factor in column K was
# incorrectly set to 1.0565 on several quote line rows. Reset those rows
# back to 1 (no surplus) to correct the pricing calculations.
$ws.Range("K16").Value = 1
$ws.Range("K17").Value = 1
$ws.Range("K21").Value = 1
$ws.Range("K24").Value = 1
$ws.Range("K27").Value = 1
$ws.Range("K28").Value = 1
$ws.Range("K32").Value = 1

# Added support for longer quotes / updated active selection on the sheet.
$ws.Activate()
$ws.Range("E3").Select()
